$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Hcrt"
$ws.Cells.Item(2,3).Value = "Hcrtr1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.1263503333333333
$ws.Cells.Item(2,8).Value = 0.379051
$ws.Cells.Item(2,9).Value = 0.270579862429723
$ws.Cells.Item(2,10).Value = 0.270579862429723
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 0.4760353333333334
$ws.Cells.Item(2,14).Value = 1.428106
$ws.Cells.Item(2,15).Value = 0.6513037878622039
$ws.Cells.Item(2,16).Value = 0.6513037878622038
$ws.Cells.Item(2,17).Value = 0.06014722304511112
$ws.Cells.Item(2,18).Value = 0.5413250074060001
$ws.Cells.Item(2,19).Value = 0.1762296893197126
$ws.Cells.Item(2,20).Value = 0.1762296893197126
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Hcrt"
$ws.Cells.Item(3,3).Value = "Hcrtr1"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.1263503333333333
$ws.Cells.Item(3,8).Value = 0.379051
$ws.Cells.Item(3,9).Value = 0.270579862429723
$ws.Cells.Item(3,10).Value = 0.270579862429723
$ws.Cells.Item(3,11).Value = 1
$ws.Cells.Item(3,12).Value = 0.3333333333333333
$ws.Cells.Item(3,13).Value = 0.02252366666666667
$ws.Cells.Item(3,14).Value = 0.06757100000000001
$ws.Cells.Item(3,15).Value = 0.03081651379494028
$ws.Cells.Item(3,16).Value = 0.03081651379494027
$ws.Cells.Item(3,17).Value = 0.002845872791222222
$ws.Cells.Item(3,18).Value = 0.02561285512100001
$ws.Cells.Item(3,19).Value = 0.0083383280631986
$ws.Cells.Item(3,20).Value = 0.008338328063198598
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Hcrt"
$ws.Cells.Item(4,3).Value = "Hcrtr1"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.1263503333333333
$ws.Cells.Item(4,8).Value = 0.379051
$ws.Cells.Item(4,9).Value = 0.270579862429723
$ws.Cells.Item(4,10).Value = 0.270579862429723
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 0.232337
$ws.Cells.Item(4,14).Value = 0.697011
$ws.Cells.Item(4,15).Value = 0.3178796983428559
$ws.Cells.Item(4,16).Value = 0.3178796983428558
$ws.Cells.Item(4,17).Value = 0.02935585739566667
$ws.Cells.Item(4,18).Value = 0.264202716561
$ws.Cells.Item(4,19).Value = 0.08601184504681179
$ws.Cells.Item(4,20).Value = 0.08601184504681178
$ws.Cells.Item(5,1).Value = "MuSCs"
$ws.Cells.Item(5,2).Value = "Hcrt"
$ws.Cells.Item(5,3).Value = "Hcrtr1"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 2
$ws.Cells.Item(5,6).Value = 0.6666666666666666
$ws.Cells.Item(5,7).Value = 0.2143626666666667
$ws.Cells.Item(5,8).Value = 0.6430880000000001
$ws.Cells.Item(5,9).Value = 0.4590587086439706
$ws.Cells.Item(5,10).Value = 0.4590587086439706
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.4760353333333334
$ws.Cells.Item(5,14).Value = 1.428106
$ws.Cells.Item(5,15).Value = 0.6513037878622039
$ws.Cells.Item(5,16).Value = 0.6513037878622038
$ws.Cells.Item(5,17).Value = 0.1020442034808889
$ws.Cells.Item(5,18).Value = 0.9183978313280002
$ws.Cells.Item(5,19).Value = 0.2989866757909499
$ws.Cells.Item(5,20).Value = 0.2989866757909498
$ws.Cells.Item(6,1).Value = "MuSCs"
$ws.Cells.Item(6,2).Value = "Hcrt"
$ws.Cells.Item(6,3).Value = "Hcrtr1"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 2
$ws.Cells.Item(6,6).Value = 0.6666666666666666
$ws.Cells.Item(6,7).Value = 0.2143626666666667
$ws.Cells.Item(6,8).Value = 0.6430880000000001
$ws.Cells.Item(6,9).Value = 0.4590587086439706
$ws.Cells.Item(6,10).Value = 0.4590587086439706
$ws.Cells.Item(6,11).Value = 1
$ws.Cells.Item(6,12).Value = 0.3333333333333333
$ws.Cells.Item(6,13).Value = 0.02252366666666667
$ws.Cells.Item(6,14).Value = 0.06757100000000001
$ws.Cells.Item(6,15).Value = 0.03081651379494028
$ws.Cells.Item(6,16).Value = 0.03081651379494027
$ws.Cells.Item(6,17).Value = 0.004828233249777779
$ws.Cells.Item(6,18).Value = 0.04345409924800001
$ws.Cells.Item(6,19).Value = 0.01414658902761439
$ws.Cells.Item(6,20).Value = 0.01414658902761439
$ws.Cells.Item(7,1).Value = "MuSCs"
$ws.Cells.Item(7,2).Value = "Hcrt"
$ws.Cells.Item(7,3).Value = "Hcrtr1"
$ws.Cells.Item(7,4).Value = "MuSCs"
$ws.Cells.Item(7,5).Value = 2
$ws.Cells.Item(7,6).Value = 0.6666666666666666
$ws.Cells.Item(7,7).Value = 0.2143626666666667
$ws.Cells.Item(7,8).Value = 0.6430880000000001
$ws.Cells.Item(7,9).Value = 0.4590587086439706
$ws.Cells.Item(7,10).Value = 0.4590587086439706
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.232337
$ws.Cells.Item(7,14).Value = 0.697011
$ws.Cells.Item(7,15).Value = 0.3178796983428559
$ws.Cells.Item(7,16).Value = 0.3178796983428558
$ws.Cells.Item(7,17).Value = 0.04980437888533334
$ws.Cells.Item(7,18).Value = 0.4482394099680001
$ws.Cells.Item(7,19).Value = 0.1459254438254063
$ws.Cells.Item(7,20).Value = 0.1459254438254063
$ws.Cells.Item(8,1).Value = "Resolving-Mac"
$ws.Cells.Item(8,2).Value = "Hcrt"
$ws.Cells.Item(8,3).Value = "Hcrtr1"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 1
$ws.Cells.Item(8,6).Value = 0.3333333333333333
$ws.Cells.Item(8,7).Value = 0.1262483333333333
$ws.Cells.Item(8,8).Value = 0.378745
$ws.Cells.Item(8,9).Value = 0.2703614289263065
$ws.Cells.Item(8,10).Value = 0.2703614289263065
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 0.4760353333333334
$ws.Cells.Item(8,14).Value = 1.428106
$ws.Cells.Item(8,15).Value = 0.6513037878622039
$ws.Cells.Item(8,16).Value = 0.6513037878622038
$ws.Cells.Item(8,17).Value = 0.06009866744111111
$ws.Cells.Item(8,18).Value = 0.5408880069700001
$ws.Cells.Item(8,19).Value = 0.1760874227515414
$ws.Cells.Item(8,20).Value = 0.1760874227515414
$ws.Cells.Item(9,1).Value = "Resolving-Mac"
$ws.Cells.Item(9,2).Value = "Hcrt"
$ws.Cells.Item(9,3).Value = "Hcrtr1"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 1
$ws.Cells.Item(9,6).Value = 0.3333333333333333
$ws.Cells.Item(9,7).Value = 0.1262483333333333
$ws.Cells.Item(9,8).Value = 0.378745
$ws.Cells.Item(9,9).Value = 0.2703614289263065
$ws.Cells.Item(9,10).Value = 0.2703614289263065
$ws.Cells.Item(9,11).Value = 1
$ws.Cells.Item(9,12).Value = 0.3333333333333333
$ws.Cells.Item(9,13).Value = 0.02252366666666667
$ws.Cells.Item(9,14).Value = 0.06757100000000001
$ws.Cells.Item(9,15).Value = 0.03081651379494028
$ws.Cells.Item(9,16).Value = 0.03081651379494027
$ws.Cells.Item(9,17).Value = 0.002843575377222222
$ws.Cells.Item(9,18).Value = 0.025592178395
$ws.Cells.Item(9,19).Value = 0.00833159670412729
$ws.Cells.Item(9,20).Value = 0.00833159670412729
$ws.Cells.Item(10,1).Value = "Resolving-Mac"
$ws.Cells.Item(10,2).Value = "Hcrt"
$ws.Cells.Item(10,3).Value = "Hcrtr1"
$ws.Cells.Item(10,4).Value = "MuSCs"
$ws.Cells.Item(10,5).Value = 1
$ws.Cells.Item(10,6).Value = 0.3333333333333333
$ws.Cells.Item(10,7).Value = 0.1262483333333333
$ws.Cells.Item(10,8).Value = 0.378745
$ws.Cells.Item(10,9).Value = 0.2703614289263065
$ws.Cells.Item(10,10).Value = 0.2703614289263065
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 0.232337
$ws.Cells.Item(10,14).Value = 0.697011
$ws.Cells.Item(10,15).Value = 0.3178796983428559
$ws.Cells.Item(10,16).Value = 0.3178796983428558
$ws.Cells.Item(10,17).Value = 0.02933215902166666
$ws.Cells.Item(10,18).Value = 0.263989431195
$ws.Cells.Item(10,19).Value = 0.08594240947063778
$ws.Cells.Item(10,20).Value = 0.08594240947063779
Write-Host "Done"
